# Update "想去人数" (F column) figures to the values captured in the latest
# gh-pages data refresh (commit 456a3b4). Applies per-sheet, per-row updates.

$wb = $excel.ActiveWorkbook

# 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1574  # was 1571
$ws.Range("F4").Value = 823  # was 819
$ws.Range("F7").Value = 1131  # was 1128
$ws.Range("F8").Value = 730  # was 729
$ws.Range("F9").Value = 780  # was 776
$ws.Range("F10").Value = 1413  # was 1400
$ws.Range("F11").Value = 279  # was 277
$ws.Range("F12").Value = 1028  # was 1026
$ws.Range("F16").Value = 47  # was 45
$ws.Range("F17").Value = 447  # was 441
$ws.Range("F18").Value = 15  # was 14
$ws.Range("F21").Value = 547  # was 545
$ws.Range("F23").Value = 752  # was 751
$ws.Range("F24").Value = 239  # was 237
$ws.Range("F25").Value = 172  # was 171
$ws.Range("F26").Value = 367  # was 366

# 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 996  # was 995
$ws.Range("F7").Value = 143  # was 141
$ws.Range("F9").Value = 584  # was 583
$ws.Range("F10").Value = 77  # was 76

# 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 218  # was 216

# 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 218  # was 216
$ws.Range("F4").Value = 1574  # was 1571
$ws.Range("F6").Value = 823  # was 819
$ws.Range("F8").Value = 996  # was 995
$ws.Range("F10").Value = 1131  # was 1128
$ws.Range("F11").Value = 730  # was 729
$ws.Range("F12").Value = 780  # was 776
$ws.Range("F13").Value = 1413  # was 1400
$ws.Range("F14").Value = 279  # was 277
$ws.Range("F15").Value = 1028  # was 1026
$ws.Range("F19").Value = 47  # was 45
$ws.Range("F20").Value = 447  # was 441
$ws.Range("F21").Value = 15  # was 14
$ws.Range("F27").Value = 143  # was 141
$ws.Range("F28").Value = 143  # was 141
$ws.Range("F29").Value = 547  # was 545
$ws.Range("F31").Value = 752  # was 751
$ws.Range("F32").Value = 239  # was 237
$ws.Range("F34").Value = 172  # was 171
$ws.Range("F35").Value = 584  # was 583
$ws.Range("F36").Value = 77  # was 76
$ws.Range("F37").Value = 77  # was 76
$ws.Range("F39").Value = 367  # was 366

